# C5-PowerPoint.pptx edit — Thu, May 07, 2020 10:07:35 AM
#
# 1) The table on slide 6 switches from the custom "Table_0" table style
#    to the built-in PowerPoint table style {27D768DD-6160-4391-88E4-E665905AB8A8}.
# 2) The deck's theme palette is restored from the custom "Integral" colour
#    scheme back to the default Office colour scheme (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), which is what actually drives the slides' look.

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 6 -------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{27D768DD-6160-4391-88E4-E665905AB8A8}")
    }
}

# --- 2) Restore the default Office theme colours ---------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
